# Update the Diet_HCOxxx model table to include PdeltaAIC as a covariate for CG as well.
# Only the numeric Chi2, p-value and (two) Estimate entries change; everything else
# (Parameter/Estimate labels, SE, DF columns, headers) stays the same.
#
# The table (sheet1) layout is:
#   Row1: Parameter | Estimate | SE | Chi2 | DF | p value   (header)
#   Row2: intrcpt             |  0.030 | 0.021 | 1.408 | 0 | 0.1590
#   Row3: Diet_HCOherbivore   | -0.015 | 0.021 | 2.631 | 2 | 0.2683
#   Row4: Diet_HCOomnivore    |  0.037 | 0.031 | 2.631 | 2 | 0.2683
#   Row5: Pvalue              | -0.036 | 0.039 | 0.857 | 1 | 0.3546
#
# All of the data cells are stored as *text* (not numbers) so that the exact
# formatting (fixed decimals, leading space for positive values, etc.) is kept.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F ("p value") cells already carry an explicit Text number format, so we can
# just replace their characters directly - this keeps their cell style untouched.
$ws.Range("F2").Characters().Text = "0.1616"
$ws.Range("F3").Characters().Text = "0.2679"
$ws.Range("F4").Characters().Text = "0.2679"
$ws.Range("F5").Characters().Text = "0.3741"

# Columns B ("Estimate") and D ("Chi2") use the General number format while still
# holding text values. Assigning a leading apostrophe forces Excel to store the
# value as text (preserving trailing zeros / leading spaces) instead of
# reinterpreting it as a number.
$ws.Range("D2").Value = "'1.400"
$ws.Range("D3").Value = "'2.634"
$ws.Range("D4").Value = "'2.634"
$ws.Range("D5").Value = "'0.790"

$ws.Range("B4").Value = "' 0.036"
$ws.Range("B5").Value = "'-0.034"
